$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while preserving it as TEXT (matching the
# original inline-string cell type), even when the new value happens to
# look like a plain number (e.g. "9.58"). For such "number-looking"
# values we temporarily force a Text number format so Excel doesn't
# silently convert the cell to a numeric type, then restore the default
# ("Normal") style so no stray formatting is left behind.
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    if ($value -match '^[0-9]+(\.[0-9]+)?$') {
        $c.NumberFormat = "@"
        $c.Value = $value
        $c.Style = "Normal"
    } else {
        $c.Value = $value
    }
}

Set-TextValue 'D2' '37.110.71'
Set-TextValue 'E2' '  +0.13%  '
Set-TextValue 'D3' '2.052.32'
Set-TextValue 'E3' '  -0.52%  '
Set-TextValue 'E4' '  +0.31%  '
Set-TextValue 'D5' '249.25'
Set-TextValue 'E5' '  -0.09%  '
Set-TextValue 'D6' '0.671'
Set-TextValue 'E6' '  -0.21%  '
Set-TextValue 'D7' '58.97'
Set-TextValue 'E7' '  +7.19%  '
Set-TextValue 'E8' '  +0.07%  '
Set-TextValue 'D9' '0.389'
Set-TextValue 'E9' '  +0.97%  '
Set-TextValue 'D10' '0.0791'
Set-TextValue 'E10' '  +0.59%  '
Set-TextValue 'E11' '  +1.93%  '
Set-TextValue 'D12' '15.96'
Set-TextValue 'E12' '  +5.80%  '
Set-TextValue 'D13' '2.354.89'
Set-TextValue 'E13' '  -0.35%  '
Set-TextValue 'D14' '0.834'
Set-TextValue 'E14' '  +1.88%  '
Set-TextValue 'D15' '5.72'
Set-TextValue 'E15' '  +7.41%  '
Set-TextValue 'D16' '2.060.58'
Set-TextValue 'E16' '  -0.16%  '
Set-TextValue 'D17' '18.76'
Set-TextValue 'E17' '  +31.47%  '
Set-TextValue 'D18' '37.053.59'
Set-TextValue 'E18' '  +0.07%  '
Set-TextValue 'D19' '75.46'
Set-TextValue 'E19' '  +2.70%  '
Set-TextValue 'D20' '0.0₃0905'
Set-TextValue 'E20' '  -2.78%  '
Set-TextValue 'D21' '5.44'
Set-TextValue 'E21' '  +1.19%  '
Set-TextValue 'D22' '238.15'
Set-TextValue 'E22' '  +0.26%  '
Set-TextValue 'E23' '  +0.00%  '
Set-TextValue 'D24' '2.42'
Set-TextValue 'E24' '  -0.73%  '
Set-TextValue 'D25' '2.22'
Set-TextValue 'E25' '  +11.74%  '
Set-TextValue 'D26' '9.58'
Set-TextValue 'E26' '  +5.91%  '
Set-TextValue 'D27' '169.22'
Set-TextValue 'E27' '  -0.61%  '
Set-TextValue 'D28' '20.13'
Set-TextValue 'E28' '  -0.24%  '
Set-TextValue 'D29' '0.126'
Set-TextValue 'E29' '  +0.82%  '
Set-TextValue 'E30' '  +6.17%  '
Set-TextValue 'D31' '4.80'
Set-TextValue 'E31' '  +3.99%  '
Set-TextValue 'D32' '0.0629'
Set-TextValue 'E32' '  +0.03%  '
Set-TextValue 'D33' '4.52'
Set-TextValue 'E33' '  +2.58%  '
Set-TextValue 'D34' '0.0900'
Set-TextValue 'E34' '  +1.06%  '
Set-TextValue 'E35' '  +0.17%  '
Set-TextValue 'D36' '2.22'
Set-TextValue 'E36' '  -3.05%  '
Set-TextValue 'D37' '1.73'
Set-TextValue 'E37' '  -2.05%  '
Set-TextValue 'E38' '  +4.58%  '
Set-TextValue 'E39' '  -1.06%  '
Set-TextValue 'E40' '  +10.84%  '
Set-TextValue 'D41' '5.11'
Set-TextValue 'E41' '  +23.78%  '
Set-TextValue 'D42' '17.69'
Set-TextValue 'E42' '  -0.25%  '
Set-TextValue 'D43' '0.0224'
Set-TextValue 'E43' '  -0.58%  '
Set-TextValue 'D44' '1.14'
Set-TextValue 'E44' '  -0.65%  '
Set-TextValue 'D45' '97.26'
Set-TextValue 'E45' '  +0.28%  '
Set-TextValue 'D46' '2.50'
Set-TextValue 'E46' '  +4.14%  '
Set-TextValue 'D47' '1.293.61'
Set-TextValue 'E47' '  -0.39%  '
Set-TextValue 'D48' '3.82'
Set-TextValue 'E48' '  -6.22%  '
Set-TextValue 'E49' '  -1.39%  '
Set-TextValue 'D50' '6.84'
Set-TextValue 'E50' '  -0.75%  '
Set-TextValue 'D51' '2.237.77'
Set-TextValue 'E51' '  -0.54%  '
